$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (A8, C9) ---
$ws.Range("A8").Value2 = "Volume 32   Number  26"
$ws.Range("C9").Value2 = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Style anchor cells (unchanged by this edit) ---
$styleGeneral = $ws.Cells.Item(33,3)   # C33 s=13 (General/text)
$styleNum166  = $ws.Cells.Item(31,10)  # J31 s=14 (#,##0)
$styleNum167  = $ws.Cells.Item(31,11)  # K31 s=15 (#,##0.0)

# Scratch cell used to force a numeric-looking string to be stored as TEXT
$scratch = $ws.Cells.Item(300,300)

function Set-TextZero($cell) {
    $scratch.NumberFormat = "@"
    $scratch.Value2 = "0"
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $styleGeneral.Copy()
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $scratch.Clear()
}

# --- Row 15 ---
$styleNum166.Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(15,3).Value2 = 1
$ws.Cells.Item(15,6).Value2 = 2
$ws.Cells.Item(15,9).Value2 = 4
$ws.Cells.Item(15,11).Value2 = 33.333333333333
$ws.Cells.Item(15,12).Value2 = 300
$ws.Cells.Item(15,13).Value2 = 100
$ws.Cells.Item(15,14).Value2 = 33.333333333333

# --- Row 16 ---
Set-TextZero ($ws.Cells.Item(16,3))
$ws.Cells.Item(16,6).Value2 = 3
$ws.Cells.Item(16,8).Value2 = -25
$ws.Cells.Item(16,9).Value2 = 12
$ws.Cells.Item(16,11).Value2 = -14.285714285714
$ws.Cells.Item(16,12).Value2 = 33.333333333333
$ws.Cells.Item(16,13).Value2 = -7.692307692307
$ws.Cells.Item(16,14).Value2 = -47.826086956521

# --- Row 17 ---
$ws.Cells.Item(17,6).Value2 = 10
$ws.Cells.Item(17,7).Value2 = 8
$ws.Cells.Item(17,8).Value2 = 25
$ws.Cells.Item(17,9).Value2 = 52
$ws.Cells.Item(17,10).Value2 = 49
$ws.Cells.Item(17,11).Value2 = 6.122448979591
$ws.Cells.Item(17,12).Value2 = 15.555555555555
$ws.Cells.Item(17,13).Value2 = 147.619047619048
$ws.Cells.Item(17,14).Value2 = -5.454545454545

# --- Row 18 ---
Set-TextZero ($ws.Cells.Item(18,3))
$ws.Cells.Item(18,7).Value2 = 1
$ws.Cells.Item(18,8).Value2 = 200
$ws.Cells.Item(18,12).Value2 = 3.333333333333
$ws.Cells.Item(18,13).Value2 = -46.551724137931
$ws.Cells.Item(18,14).Value2 = -78.321678321678

# --- Row 19 ---
$ws.Cells.Item(19,3).Value2 = 5
$ws.Cells.Item(19,4).Value2 = 7
$ws.Cells.Item(19,5).Value2 = -28.571428571428
$ws.Cells.Item(19,6).Value2 = 23
$ws.Cells.Item(19,7).Value2 = 35
$ws.Cells.Item(19,8).Value2 = -34.285714285714
$ws.Cells.Item(19,9).Value2 = 127
$ws.Cells.Item(19,10).Value2 = 151
$ws.Cells.Item(19,11).Value2 = -15.894039735099
$ws.Cells.Item(19,12).Value2 = -11.805555555555
$ws.Cells.Item(19,13).Value2 = 98.4375
$ws.Cells.Item(19,14).Value2 = 41.111111111111

# --- Row 20 ---
Set-TextZero ($ws.Cells.Item(20,3))
$styleNum166.Copy()
$ws.Cells.Item(20,4).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(20,4).Value2 = 4
$styleNum167.Copy()
$ws.Cells.Item(20,5).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(20,5).Value2 = -100
$ws.Cells.Item(20,6).Value2 = 1
$ws.Cells.Item(20,7).Value2 = 6
$ws.Cells.Item(20,8).Value2 = -83.333333333333
$ws.Cells.Item(20,10).Value2 = 23
$ws.Cells.Item(20,11).Value2 = -43.478260869565
$ws.Cells.Item(20,13).Value2 = -23.529411764705
$ws.Cells.Item(20,14).Value2 = -96.306818181818

# --- Row 21 ---
$ws.Cells.Item(21,3).Value2 = 8
$ws.Cells.Item(21,4).Value2 = 13
$ws.Cells.Item(21,5).Value2 = -38.461538461538
$ws.Cells.Item(21,6).Value2 = 42
$ws.Cells.Item(21,7).Value2 = 54
$ws.Cells.Item(21,8).Value2 = -22.222222222222
$ws.Cells.Item(21,9).Value2 = 239
$ws.Cells.Item(21,10).Value2 = 257
$ws.Cells.Item(21,11).Value2 = -7.003891050583
$ws.Cells.Item(21,12).Value2 = -11.481481481481
$ws.Cells.Item(21,13).Value2 = 36.571428571428
$ws.Cells.Item(21,14).Value2 = -64.167916041979

# --- Row 24 ---
$ws.Cells.Item(24,3).Value2 = 11
$ws.Cells.Item(24,4).Value2 = 7
$ws.Cells.Item(24,5).Value2 = 57.142857142857
$ws.Cells.Item(24,6).Value2 = 32
$ws.Cells.Item(24,8).Value2 = 6.666666666666
$ws.Cells.Item(24,9).Value2 = 221
$ws.Cells.Item(24,10).Value2 = 199
$ws.Cells.Item(24,11).Value2 = 11.055276381909
$ws.Cells.Item(24,12).Value2 = -7.142857142857
$ws.Cells.Item(24,13).Value2 = -10.526315789473

# --- Row 25 ---
$ws.Cells.Item(25,3).Value2 = 6
$ws.Cells.Item(25,4).Value2 = 1
$ws.Cells.Item(25,5).Value2 = 500
$ws.Cells.Item(25,6).Value2 = 18
$ws.Cells.Item(25,7).Value2 = 13
$ws.Cells.Item(25,8).Value2 = 38.461538461538
$ws.Cells.Item(25,9).Value2 = 125
$ws.Cells.Item(25,10).Value2 = 95
$ws.Cells.Item(25,11).Value2 = 31.578947368421
$ws.Cells.Item(25,12).Value2 = 12.612612612612

# --- Row 26 ---
$ws.Cells.Item(26,3).Value2 = 8
$ws.Cells.Item(26,4).Value2 = 5
$ws.Cells.Item(26,5).Value2 = 60
$ws.Cells.Item(26,6).Value2 = 20
$ws.Cells.Item(26,7).Value2 = 15
$ws.Cells.Item(26,8).Value2 = 33.333333333333
$ws.Cells.Item(26,9).Value2 = 111
$ws.Cells.Item(26,10).Value2 = 82
$ws.Cells.Item(26,11).Value2 = 35.365853658536
$ws.Cells.Item(26,12).Value2 = 16.842105263157
$ws.Cells.Item(26,13).Value2 = 6.730769230769

# --- Row 27 ---
$styleNum166.Copy()
$ws.Cells.Item(27,3).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(27,3).Value2 = 1
$ws.Cells.Item(27,6).Value2 = 2
$ws.Cells.Item(27,9).Value2 = 4
$ws.Cells.Item(27,11).Value2 = -42.857142857142
$ws.Cells.Item(27,12).Value2 = 100

# --- Row 28 ---
Set-TextZero ($ws.Cells.Item(28,3))
$ws.Cells.Item(28,5).Value2 = -100
$ws.Cells.Item(28,7).Value2 = 4
$ws.Cells.Item(28,8).Value2 = -75
$ws.Cells.Item(28,10).Value2 = 7
$ws.Cells.Item(28,11).Value2 = 14.285714285714
